$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.533.57'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.825.54'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'315.54"
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = "'0.5115"
$ws.Range("E7").Value = '  -5.52%  '
$ws.Range("D8").Value = "'0.3955"
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("D9").Value = "'0.08213"
$ws.Range("E9").Value = '  +6.43%  '
$ws.Range("E10").Value = '  -0.84%  '
$ws.Range("D11").Value = "'41.70"
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("D12").Value = "'21.18"
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = "'6.344"
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("E14").Value = '  +0.10%  '
$ws.Range("D15").Value = "'7.552"
$ws.Range("E15").Value = '  -1.21%  '
$ws.Range("D16").Value = '1.823.95'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = "'0.00001127"
$ws.Range("E17").Value = '  +3.36%  '
$ws.Range("D18").Value = "'92.91"
$ws.Range("E18").Value = '  +3.11%  '
$ws.Range("D19").Value = "'0.06662"
$ws.Range("E19").Value = '  +1.04%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").Value = "'6.097"
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").Value = '28.576.67'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("D24").Value = "'11.43"
$ws.Range("E24").Value = '  +2.08%  '
$ws.Range("D25").Value = "'2.261"
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("D26").Value = "'21.40"
$ws.Range("E26").Value = '  +2.85%  '
$ws.Range("D27").Value = "'156.74"
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("D28").Value = '2.036.75'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = "'2.409"
$ws.Range("E29").Value = '  -2.23%  '
$ws.Range("D30").Value = "'127.16"
$ws.Range("E30").Value = '  +2.23%  '
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("D32").Value = "'0.1090"
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("D33").Value = "'5.764"
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").Value = "'3.656"
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").Value = "'0.07075"
$ws.Range("E35").Value = '  -6.79%  '
$ws.Range("E36").Value = '  -0.91%  '
$ws.Range("D37").Value = "'5.286"
$ws.Range("E37").Value = '  +1.28%  '
$ws.Range("D38").Value = "'0.02353"
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").Value = "'8.812"
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("D40").Value = "'0.6330"
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("D41").Value = "'11.29"
$ws.Range("E41").Value = '  -0.79%  '
$ws.Range("D42").Value = "'1.184"
$ws.Range("E42").Value = '  -0.50%  '
$ws.Range("D43").Value = "'1.398"
$ws.Range("E43").Value = '  -0.24%  '
$ws.Range("D44").Value = "'13.54"
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").Value = "'0.5945"
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").Value = "'125.35"
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("D49").Value = "'1.194"
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").Value = "'0.06939"
$ws.Range("E50").Value = '  +0.31%  '
$ws.Range("D51").Value = "'1.084"
$ws.Range("E51").Value = '  +4.38%  '
